$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.06328177979961902
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 16.98373111632243
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 23.86602580940113

$ws.Range("B3").Value = 0.1554434735375247
$ws.Range("C3").Value = 0.05231270169004087
$ws.Range("D3").Value = 0.7127328510149897
$ws.Range("E3").Value = 6.48142807727062
$ws.Range("G3").Value = 7.401917103513176

$ws.Range("B4").Value = 0.006876353814593728
$ws.Range("C4").Value = 0.004309184025731883
$ws.Range("D4").Value = 3.082599426703578
$ws.Range("E4").Value = 246.9852506941017
$ws.Range("G4").Value = 250.0790356586456

$ws.Range("B5").Value = 0.7287194209349384
$ws.Range("C5").Value = 1.65323645889881
$ws.Range("D5").Value = 0.1529057820181812
$ws.Range("E5").Value = 0.4998867070740569
$ws.Range("G5").Value = 3.034748368925986
